# Apply updated coin price/volume data scraped on 2023-02-15 (05:40 UTC run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking / percentage text must be written with a leading apostrophe so
# Excel stores it as literal text (matching the sheet's existing text-formatted cells)
# instead of silently converting it to a Number/Percentage value.
$q = "'"

# Row 2
$ws.Range("D2").Value = $q + '297.30'
$ws.Range("E2").Value = $q + '1.67%'

# Row 3
$ws.Range("D3").Value = $q + '41.88'
$ws.Range("E3").Value = $q + '3.57%'

# Row 4
$ws.Range("D4").Value = $q + '5.004'
$ws.Range("E4").Value = $q + '-0.46%'

# Row 5
$ws.Range("D5").Value = $q + '0.07519'
$ws.Range("E5").Value = $q + '2.61%'

# Row 6
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").Value = $q + '1.584'
$ws.Range("E6").Value = $q + '3.32%'

# Row 7
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = $q + '0.9257'
$ws.Range("E7").Value = $q + '-0.17%'

# Row 8
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = $q + '2.401'
$ws.Range("E8").Value = $q + '1.36%'

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = $q + '0.1194'
$ws.Range("E9").Value = $q + '0.69%'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = $q + '0.1823'
$ws.Range("E10").Value = $q + '4.48%'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = $q + '0.08916'
$ws.Range("E11").Value = $q + '3.14%'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = $q + '0.04077'
$ws.Range("E12").Value = $q + '-5.79%'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = $q + '0.1049'
$ws.Range("E13").Value = $q + '-0.54%'

# Row 14
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D14").Value = $q + '0.005788'
$ws.Range("E14").Value = $q + '-3.07%'

# Row 15
$ws.Range("B15").Value = 'LEO'
$ws.Range("C15").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D15").Value = $q + '3.357'
$ws.Range("E15").Value = $q + '0.53%'

# Row 16
$ws.Range("B16").Value = 'GateToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D16").Value = $q + '4.374'
$ws.Range("E16").Value = $q + '1.73%'

# Row 17
$ws.Range("B17").Value = 'BitpandaEcosystemToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D17").Value = $q + '0.3313'
$ws.Range("E17").Value = $q + '0.72%'

# Row 18
$ws.Range("B18").Value = 'MCDex'
$ws.Range("C18").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D18").Value = $q + '8.118'
$ws.Range("E18").Value = $q + '1.84%'

# Row 19
$ws.Range("B19").Value = 'ProBitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D19").Value = $q + '0.1391'
$ws.Range("E19").Value = $q + '0.03%'

# Row 20
$ws.Range("B20").Value = 'ZBToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D20").Value = $q + '0.3102'
$ws.Range("E20").Value = $q + '11.06%'

# Row 21
$ws.Range("B21").Value = 'BitForexToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D21").Value = $q + '0.001283'
$ws.Range("E21").Value = $q + '0.68%'

# Row 22
$ws.Range("D22").Value = $q + '0.04076'
$ws.Range("E22").Value = $q + '3.68%'

# Row 23
$ws.Range("D23").Value = $q + '0.001266'
$ws.Range("E23").Value = $q + '0.34%'

# Row 24
$ws.Range("D24").Value = $q + '0.003903'
$ws.Range("E24").Value = $q + '3.25%'

# Row 25
$ws.Range("E25").Value = $q + '-3.94%'

# Row 38
$ws.Range("D38").Value = $q + '0.02406'
$ws.Range("E38").Value = $q + '5.76%'

# Row 39
$ws.Range("D39").Value = $q + '0.05211'
$ws.Range("E39").Value = $q + '4.87%'

# Row 40
$ws.Range("D40").Value = $q + '0.006305'
$ws.Range("E40").Value = $q + '11.92%'

# Row 41
$ws.Range("D41").Value = $q + '0.007834'
$ws.Range("E41").Value = $q + '1.85%'

# Row 42
$ws.Range("E42").Value = $q + '3.13%'

# Row 43
$ws.Range("D43").Value = $q + '0.007409'
$ws.Range("E43").Value = $q + '0.60%'

# Row 44
$ws.Range("D44").Value = $q + '0.007262'
$ws.Range("E44").Value = $q + '-0.56%'

# Row 45
$ws.Range("D45").Value = $q + '0.2968'
$ws.Range("E45").Value = $q + '1.63%'

# Row 46
$ws.Range("D46").Value = $q + '0.00006596'
$ws.Range("E46").Value = $q + '4.59%'

# Row 47
$ws.Range("E47").Value = $q + '-0.03%'

# Row 48
$ws.Range("E48").Value = $q + '48.36%'

# Row 49
$ws.Range("D49").Value = $q + '0.004203'
$ws.Range("E49").Value = $q + '0.04%'

# Row 50
$ws.Range("E50").Value = $q + '-0.03%'

# Row 51
$ws.Range("D51").Value = $q + '0.0002001'
$ws.Range("E51").Value = $q + '-0.03%'
